$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so multi-part version-like
# numbers (e.g. "42.121.02") and plain decimals alike stay literal strings
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.121.02"
$ws.Range("E2").Value = "  -3.65%  "
$ws.Range("D3").Value = "2.203.47"
$ws.Range("E3").Value = "  -3.33%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "106.12"
$ws.Range("E5").Value = "  -14.15%  "
$ws.Range("D6").Value = "294.57"
$ws.Range("E6").Value = "  +10.48%  "
$ws.Range("D7").Value = "0.618"
$ws.Range("E7").Value = "  -3.52%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").Value = "  -5.60%  "
$ws.Range("D10").Value = "43.21"
$ws.Range("E10").Value = "  -10.24%  "
$ws.Range("D11").Value = "0.0902"
$ws.Range("E11").Value = "  -4.73%  "
$ws.Range("D12").Value = "53.94"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").Value = "8.64"
$ws.Range("E13").Value = "  -7.14%  "
$ws.Range("E14").Value = "  -3.57%  "
$ws.Range("D15").Value = "0.932"
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("E16").Value = "  -4.31%  "
$ws.Range("D17").Value = "2.537.41"
$ws.Range("E17").Value = "  -3.31%  "
$ws.Range("D18").Value = "2.225.84"
$ws.Range("E18").Value = "  -2.44%  "
$ws.Range("D19").Value = "42.079.91"
$ws.Range("E19").Value = "  -3.52%  "
$ws.Range("D20").Value = "7.24"
$ws.Range("E20").Value = "  +3.68%  "
$ws.Range("E21").Value = "  -5.79%  "
$ws.Range("D22").Value = "72.01"
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").Value = "3.44"
$ws.Range("E23").Value = "  +18.89%  "
$ws.Range("D24").Value = "2.25"
$ws.Range("E24").Value = "  -7.66%  "
$ws.Range("D25").Value = "226.17"
$ws.Range("E25").Value = "  -4.10%  "
$ws.Range("D26").Value = "8.83"
$ws.Range("E26").Value = "  -7.21%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("D28").Value = "11.45"
$ws.Range("E28").Value = "  -4.18%  "
$ws.Range("D29").Value = "3.95"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").Value = "2.24"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").Value = "37.55"
$ws.Range("E31").Value = "  -12.05%  "
$ws.Range("D32").Value = "3.19"
$ws.Range("E32").Value = "  -5.23%  "
$ws.Range("D33").Value = "172.64"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").Value = "20.70"
$ws.Range("E34").Value = "  -4.64%  "
$ws.Range("D35").Value = "0.0867"
$ws.Range("E35").Value = "  -5.90%  "
$ws.Range("D36").Value = "4.96"
$ws.Range("E36").Value = "  +6.86%  "
$ws.Range("D37").Value = "5.44"
$ws.Range("E37").Value = "  -5.83%  "
$ws.Range("D38").Value = "4.23"
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("D39").Value = "0.124"
$ws.Range("E39").Value = "  -4.14%  "
$ws.Range("D40").Value = "0.0358"
$ws.Range("E40").Value = "  -5.32%  "
$ws.Range("E41").Value = "  -5.57%  "
$ws.Range("D42").Value = "2.43"
$ws.Range("E42").Value = "  -4.58%  "
$ws.Range("D43").Value = "69.77"
$ws.Range("E43").Value = "  -6.00%  "
$ws.Range("D44").Value = "0.226"
$ws.Range("E44").Value = "  -5.30%  "
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").Value = "12.49"
$ws.Range("E46").Value = "  -10.54%  "
$ws.Range("E47").Value = "  -6.76%  "
$ws.Range("D48").Value = "5.36"
$ws.Range("E48").Value = "  -4.85%  "
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("D50").Value = "101.79"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").Value = "8.38"
$ws.Range("E51").Value = "  -2.69%  "

# Restore the original (default) cell style now that the text values are set,
# so no stray number-format/style is left attached to the cells.
$ws.Range("D2:D51").Style = "Normal"

